$d = $word.ActiveDocument

$replacements = @(
    @("820÷3=", "430÷8="),
    @("713÷5=", "152÷3="),
    @("586÷7=", "314÷6="),
    @("826÷9=", "715÷2="),
    @("526÷2=", "442÷6="),
    @("927÷4=", "378÷7="),
    @("444÷8=", "931÷9="),
    @("129÷4=", "981÷6="),
    @("228÷7=", "684÷8="),
    @("806÷2=", "464÷7="),
    @("698÷4=", "573÷7="),
    @("888÷4=", "554÷4="),
    @("688÷6=", "497÷9="),
    @("126÷6=", "930÷2="),
    @("454÷7=", "643÷3="),
    @("663÷5=", "745÷6="),
    @("845÷2=", "495÷3="),
    @("127÷8=", "404÷7="),
    @("249÷5=", "900÷6="),
    @("999÷6=", "267÷6="),
    @("228÷6=", "581÷9="),
    @("836÷5=", "489÷9="),
    @("317÷4=", "691÷6="),
    @("270÷5=", "547÷6="),
    @("185÷2=", "702÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
